# Fruta / hortaliza, semanal
# Update the weekly price records (rows 2-13) for Breva at
# Vega Central Mapocho de Santiago with the latest reported values.
# Row 6 is unchanged; only the date/quality/price/origin columns of
# the other rows are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44187
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("R2").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S2").Value = 2000

# Row 3
$ws.Range("D3").Value = 44187
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 1714

# Row 4
$ws.Range("D4").Value = 44553
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("Q4").Value = "$/bandeja 6 kilos"
$ws.Range("S4").Value = 3667
$ws.Range("T4").Value = 6

# Row 5
$ws.Range("D5").Value = 44553
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = "$/bandeja 6 kilos"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 6

# Row 6 unchanged

# Row 7
$ws.Range("D7").Value = 44204
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 110
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7500
$ws.Range("P7").Value = 7318
$ws.Range("Q7").Value = "$/bandeja 7 kilos"
$ws.Range("S7").Value = 1045
$ws.Range("T7").Value = 7

# Row 8
$ws.Range("D8").Value = 44572
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 3333

# Row 9
$ws.Range("D9").Value = 44558
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 22000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 22000
$ws.Range("Q9").Value = "$/bandeja 6 kilos"
$ws.Range("S9").Value = 3667
$ws.Range("T9").Value = 6

# Row 10
$ws.Range("D10").Value = 44558
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = "$/bandeja 6 kilos"
$ws.Range("S10").Value = 3000
$ws.Range("T10").Value = 6

# Row 11
$ws.Range("D11").Value = 44550
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("Q11").Value = "$/bandeja 7 kilos"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 3429
$ws.Range("T11").Value = 7

# Row 12
$ws.Range("D12").Value = 44189
$ws.Range("L12").Value = "Especial"
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = "$/bandeja 7 kilos"
$ws.Range("S12").Value = 2143
$ws.Range("T12").Value = 7

# Row 13
$ws.Range("D13").Value = 44189
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 13000
$ws.Range("Q13").Value = "$/bandeja 7 kilos"
$ws.Range("R13").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S13").Value = 1857
$ws.Range("T13").Value = 7
